$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.549410343170166
$ws.Range("B1").Value = 1.421885132789612
$ws.Range("C1").Value = 4.634285449981689
$ws.Range("D1").Value = 1.930008769035339
$ws.Range("E1").Value = 0.6884947419166565
